# Daily attendance processing - 2025-10-17 12:36:14
# Reverses the order of the comma-separated "Recorded By" names/emails
# in column G for every row on the active sheet that has more than one
# entry (i.e. contains a comma). Single-value cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Text

    if ($txt -like "*,*") {
        $parts = $txt -split ", "

        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }

        $newVal = $reversed -join ", "
        $cell.Value = $newVal
    }
}
